$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = 112182711
$ws.Range("B12").Value = 98535
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 222498
$ws.Range("F12").Value = 'Blåsippa'
$ws.Range("G12").Value = 'Hepatica nobilis'
$ws.Range("H12").Value = 'Schreb.'
$ws.Range("I12").Value = "'1"
$ws.Range("I12").Style = "Normal"
$ws.Range("P12").Value = 'Stjärnsund, Dlr'
$ws.Range("Q12").Value = 565668.8439373589
$ws.Range("R12").Value = 6699889.440735213
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = 'Dalarna'
$ws.Range("U12").Value = 'Hedemora'
$ws.Range("V12").Value = 'Dalarna'
$ws.Range("W12").Value = 'Husby'
$ws.Range("Y12").Value = "'2023-09-05"
$ws.Range("Y12").Style = "Normal"
$ws.Range("Z12").Value = '00:00'
$ws.Range("AA12").Value = "'2023-09-05"
$ws.Range("AA12").Style = "Normal"
$ws.Range("AB12").Value = '00:00'
$ws.Range("AC12").Value = 'Påträffad under Sveaskogs naturvärdesinventering'
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
# AT12 is blank/empty in source (left unset; reads as empty either way)
$ws.Range("AW12").Value = 'Mimmi Persson'
$ws.Range("AX12").Value = 'Mimmi Persson'
# AY12 is blank/empty in source (left unset; reads as empty either way)

# Row 13
$ws.Range("A13").Value = 112183915
$ws.Range("B13").Value = 89405
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = 'Ullticka'
$ws.Range("G13").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H13").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I13").Value = "'1"
$ws.Range("I13").Style = "Normal"
$ws.Range("P13").Value = 'Stjärnsund, Dlr'
$ws.Range("Q13").Value = 565683.2851149164
$ws.Range("R13").Value = 6699911.431222402
$ws.Range("S13").Value = 5
$ws.Range("T13").Value = 'Dalarna'
$ws.Range("U13").Value = 'Hedemora'
$ws.Range("V13").Value = 'Dalarna'
$ws.Range("W13").Value = 'Husby'
$ws.Range("Y13").Value = "'2023-09-05"
$ws.Range("Y13").Style = "Normal"
$ws.Range("Z13").Value = '00:00'
$ws.Range("AA13").Value = "'2023-09-05"
$ws.Range("AA13").Style = "Normal"
$ws.Range("AB13").Value = '00:00'
$ws.Range("AC13").Value = 'Påträffad under Sveaskogs naturvärdesinventering'
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
# AT13 is blank/empty in source (left unset; reads as empty either way)
$ws.Range("AW13").Value = 'Mimmi Persson'
$ws.Range("AX13").Value = 'Mimmi Persson'
# AY13 is blank/empty in source (left unset; reads as empty either way)
